# Delete the "lxml" acknowledgment row from the Wordless acknowledgments sheet.
# The table (Name, Home Page, Version, Authors, License, License URL) is sorted
# alphabetically by Name; "lxml" currently occupies row 11. Removing it shifts
# every following row up by one and shrinks the used range from A1:F37 to A1:F36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "lxml" row by its Name value (column A) rather than hard-coding the
# row index, so the script is resilient if the sheet differs slightly.
$lxmlRow = $null
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value -eq "lxml") {
        $lxmlRow = $r
        break
    }
}

if ($lxmlRow -eq $null) {
    $lxmlRow = 11
}

# Remove the entire row, shifting all following rows up.
$ws.Rows.Item($lxmlRow).EntireRow.Delete()

# Reset the view back to the top-left of the frozen pane and select A9, matching
# where the edit left the selection after the row removal.
$ws.Range("A9").Select()
